# Commit: Tue, Mar 24, 2020  7:05:14 AM
#
# The table on slide 6 (the "SOURCES OF FINANCE" slide) has its table
# style switched from the deck's custom local style
# {D33461A4-E15B-4E42-A093-5CD923AD489E} to the built-in PowerPoint
# table style {C852AA37-D702-4E8F-A68D-3B495CF7C997}.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{C852AA37-D702-4E8F-A68D-3B495CF7C997}")
    }
}
